# Generate Report for Handoff
# Adds a new file entry (69f442ba-1f63-4971-9535-3cd24131af35) as row 3
# on all three worksheets (Overview, zh-cn, de-de), mirroring the existing
# 18e181b3-... row, including hyperlinks and date-styled text cells.

$wb = $excel.ActiveWorkbook

$newId   = "69f442ba-1f63-4971-9535-3cd24131af35"
$newMd   = "$newId.md"
$newTok  = "ed8a66da40c2beb4ee404ef04492513c1d33ec3f"
$zhXlf   = "$newId.$newTok.zh-cn.xlf"
$deXlf   = "$newId.$newTok.de-de.xlf"

$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/de3eb2d856f4ded5769487b38265b668987b043b/e2e/$newMd"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b7ff6fb8529c00fb632cbc701aa44b748993e99/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1241332357d20ed2b1783d324dbe974ea3a3752d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A3").Value2 = $newMd
$ws1.Hyperlinks.Add($ws1.Range("A3"), $mdUrl, "", "", $newMd)

$ws1.Range("B3").Value2 = "Ready for handoff"
$ws1.Range("C3").Value2 = "Ready for handoff"

$ws1.Range("D3").Value2 = "2016-03-22 18:37:28"
$ws1.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A3").Value2 = $newMd
$ws2.Hyperlinks.Add($ws2.Range("A3"), $mdUrl, "", "", $newMd)

$ws2.Range("B3").Value2 = ".md"
$ws2.Range("C3").Value2 = "Ready for handoff"

$ws2.Range("D3").Value2 = $zhXlf
$ws2.Hyperlinks.Add($ws2.Range("D3"), $zhXlfUrl, "", "", $zhXlf)

$ws2.Range("E3").Value2 = "2016-03-22 18:37:24"
$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range("J3").Value2 = "Include"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A3").Value2 = $newMd
$ws3.Hyperlinks.Add($ws3.Range("A3"), $mdUrl, "", "", $newMd)

$ws3.Range("B3").Value2 = ".md"
$ws3.Range("C3").Value2 = "Ready for handoff"

$ws3.Range("D3").Value2 = $deXlf
$ws3.Hyperlinks.Add($ws3.Range("D3"), $deXlfUrl, "", "", $deXlf)

$ws3.Range("E3").Value2 = "2016-03-22 18:37:28"
$ws3.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Range("J3").Value2 = "Include"

Write-Output "Applied handoff report update"
